$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add four new reference-data rows for PROPERTY_AREA_TYPE (rows 114-117) ---
# Seed formatting for the new rows by copying the format of the last existing
# data row (113), which reuses the existing style records (s=3 / s=8 / s=9)
# instead of allocating new ones.
$ws.Range("A113:N113").Copy()
$ws.Range("A114:N117").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows("114:117").RowHeight = 38.65

$rows = @(
  @{ Row = 114; Key = -1; Value = "No Data" },
  @{ Row = 115; Key = 1;  Value = "Carpet" },
  @{ Row = 116; Key = 2;  Value = "Built Up" },
  @{ Row = 117; Key = 3;  Value = "Super Built Up" }
)

foreach ($r in $rows) {
  $i = $r.Row
  $ws.Range("A$i").Value = "PROPERTY_AREA_TYPE"
  $ws.Range("B$i").Value = $r.Key
  $ws.Range("C$i").Value = $r.Value
  $ws.Range("D$i").Value = "DEFAULT"
  $ws.Range("F$i").Value = "REALTY"
  $ws.Range("G$i").Value = "PROPERTY"
  $ws.Range("H$i").Value = 1
  $ws.Range("I$i").Value = "Active"
  $ws.Range("J$i").Value = "System-User"
  $ws.Range("L$i").Value = "System-User"
  $ws.Range("N$i").Formula = "=CONCATENATE(""INSERT INTO APP_REF_DATA (DATA_TYPE, DATA_KEY, DATA_VALUE, GROUP_NAME, SUB_GROUP_NAME, APP_NAME, MODULE_NAME, STATUS_ID, STATUS_TITLE, CREATED_BY, UPDATED_BY) "",""VALUES ("",""'"",A$i,""', '"",B$i,""', '"",C$i,""', '"",D$i,""', '"",E$i,""', '"",F$i,""', '"",G$i,""', '"",H$i,""', '"",I$i,""', '"",J$i,""', '"",L$i,""');"")"
}
